$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "3 cases extra thin ravioli sheets"
$ws.Range("B8").Value = "3 RAVSHE 4E"
